$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('A3').Value = 'POLYGON Z ((-1268746.656403215 2958669.4221782363 0, -1269148.7719312687 2958739.0154197128 0, -1269547.4615757428 2958814.6065749167 0, -1269546.6619563145 2958817.763578659 0, -1269474.6916636892 2959210.0765392943 0, -1269075.0058828911 2959137.6925412486 0, -1268673.8272412254 2959065.033669607 0, -1268745.8565798467 2958672.5777471345 0, -1268746.656403215 2958669.4221782363 0))'
$ws.Range('D3').Value = 'None Listed, Dept. of Transportation'
$ws.Range('E3').Value = 'Subsurface, subsurface'
$ws.Range('E7').Value = 'Subsurface, Surface'
$ws.Range('D8').Value = 'Common Schools, Oklahoma State University, Normal Schools, Connors State College, University Preparatory, Langston University, State Educational Institutions'
$ws.Range('D15').Value = 'Dept. of Parks and Recreation, 100.00% Dept. of Parks and Recreation, Public School (Indemnity, Schools, Common Schools), 100.00% Public School (Indemnity, Schools, Common Schools), 100.00% Normal School, 100.00% University of Idaho, State Hospital South (Insane Asylum), Normal School, Agricultural College, 100.00% Dept. of Fish and Game, General Fund, 100.00% General Fund'
$ws.Range('E21').Value = 'Surface, Subsurface'
$ws.Range('E28').Value = 'Surface, Subsurface'
$ws.Range('D31').Value = 'Oklahoma State University, Common Schools, Other Agency'
$ws.Range('E31').Value = 'surface, subsurface'
$ws.Range('E40').Value = 'Subsurface, Surface'
$ws.Range('E41').Value = 'Subsurface, Surface'
$ws.Range('E46').Value = 'Subsurface, Surface'
$ws.Range('D48').Value = 'Common Schools, Oklahoma State University, Public Building, State Educational Institutions, Other Agency'
$ws.Range('E48').Value = 'subsurface, surface'
$ws.Range('E53').Value = 'Surface, Subsurface'
$ws.Range('A61').Value = 'MULTIPOLYGON Z (((-150727.67552972207 2180759.685291849 0, -151058.23785022597 2180767.194620998 0, -151050.6643970744 2181169.403240565 0, -150647.36731614405 2181162.756169655 0, -150654.51162035248 2180758.678847399 0, -150727.67552972207 2180759.685291849 0)), ((-156213.95280016185 2185771.815552367 0, -156194.81600412223 2186180.9574394897 0, -156176.03287459628 2186582.5349008576 0, -156157.59366610707 2186976.756056978 0, -156139.48908621314 2187363.8209823803 0, -155744.84600765648 2187348.445937388 0, -155348.1509325643 2187332.9886809853 0, -154949.3876756076 2187317.4484555325 0, -154548.54004262667 2187301.8247038866 0, -154552.12219585586 2186910.1019930975 0, -154555.77139012987 2186511.059859929 0, -154559.48947036255 2186104.491345304 0, -154563.2783467011 2185690.181225061 0, -154979.2966338558 2185710.759282247 0, -155393.06976225073 2185731.2236910015 0, -155804.6157940152 2185751.5754990443 0, -156213.95280016185 2185771.815552367 0)), ((-154756.49514332798 2195401.4538690597 0, -155153.7515568258 2195408.221370756 0, -155553.12051018604 2195415.229121782 0, -155950.92415796503 2195422.2057534335 0, -155947.17202052783 2195827.355199152 0, -155943.4074702889 2196233.814884058 0, -155941.47868866843 2196637.651863205 0, -155939.54341035578 2197042.7984314207 0, -155541.2189184494 2197040.569327221 0, -155353.65413790673 2197039.518888024 0, -155396.12102644477 2196817.527901084 0, -155160.51979592125 2196817.2483607503 0, -155161.41836884365 2197038.4414134705 0, -155141.31684328214 2197038.3286947077 0, -154739.82762817875 2197036.076583675 0, -154336.74182022538 2197033.8129475242 0, -154341.998905905 2196622.0189541783 0, -154347.23910693155 2196211.5658134827 0, -154352.46280488887 2195802.447130641 0, -154357.67009668285 2195394.6561956443 0, -154756.49514332798 2195401.4538690597 0)), ((-150662.53300438932 2201565.8143678103 0, -150265.73542722332 2201558.079416251 0, -150266.11107496227 2201956.7167395214 0, -149870.38817961563 2201949.4151801527 0, -149473.8411775989 2201942.0961809494 0, -149469.64870493274 2201542.554079617 0, -149465.41586004765 2201139.138079935 0, -149461.14186428962 2200731.7913708803 0, -149456.82626367526 2200320.456549403 0, -149861.1348939628 2200329.6388499835 0, -150264.58700048906 2200338.799250088 0, -150667.18486721066 2200347.937866776 0, -151068.93166755023 2200357.0549574546 0, -151065.42192813862 2200766.4810241363 0, -151061.94612094056 2201171.95410604 0, -151058.50384008256 2201573.53086171 0, -150662.53300438932 2201565.8143678103 0)), ((-144985.75394429002 2186761.1955177807 0, -144993.57760826222 2186358.8521826295 0, -145001.40205040292 2185956.507569622 0, -145399.94308498083 2185961.978390855 0, -145391.28270177462 2186365.0471574035 0, -145382.62273704054 2186768.1136125233 0, -144985.75394429002 2186761.1955177807 0)), ((-146176.35634732372 2186781.936397038 0, -146186.68968903017 2186377.4237080743 0, -146584.39309958878 2186383.6114961575 0, -146573.22349830132 2186788.8468318544 0, -146562.05367514814 2187194.0780807077 0, -146166.02163888028 2187186.444051743 0, -146176.35634732372 2186781.936397038 0)), ((-150383.6975485528 2192048.570456311 0, -149985.44773027493 2192037.957515076 0, -149604.66534376488 2192027.8079793444 0, -149617.6731232064 2191624.1751634963 0, -149623.36915654116 2191447.4289275855 0, -149642.32078511678 2191298.056048375 0, -149657.298111479 2191222.1884050462 0, -149692.36435896214 2191044.5601496864 0, -149675.98934608878 2191015.215095343 0, -149620.52489968407 2190931.2250189628 0, -149623.67058447894 2190819.8473241213 0, -149635.45273980167 2190820.0542867915 0, -149646.50565726828 2190435.131514635 0, -150028.3764674055 2190441.381687129 0, -150424.5544889028 2190447.8637507427 0, -150414.70544404458 2190833.735152485 0, -150404.3979716605 2191237.5631112135 0, -150394.06204562637 2191642.5065784375 0, -150383.6975485528 2192048.570456311 0)))'
$ws.Range('E61').Value = 'Surface, Subsurface'
$ws.Range('F61').Value = 12
$ws.Range('I61').Value = 2337.21
$ws.Range('J61').Value = 6
$ws.Range('K61').Value = 4674.42
$ws.Range('E75').Value = 'Subsurface, Surface'
$ws.Range('E77').Value = 'Surface, Subsurface'
